$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The answer paragraph to "¿Por qué se dice que un protocolo es
# independiente de la tecnología?" was split across two runs (with a
# _GoBack bookmark sitting between them). Re-assert the whole sentence as
# one continuous piece of text so it collapses back into a single run and
# drops the now-stray bookmark.
$oldAnswer = "Esto se presenta debido a que los protocolos especifican las funcionalidades de la red m" + `
    "as no como se deben cumplir, ni en qué tecnología se debe usar. Por ejemplo, los protocolos como " + `
    "HTTP, DNS o TCP se pueden usar en diferentes dispositivos independientemente del sistema operativo " + `
    "con el que cuenta cada uno."
$newAnswer = "Esto se presenta debido a que los protocolos especifican las funcionalidades de la red mas no como se deben cumplir, ni en qué tecnología se debe usar. Por ejemplo, los protocolos como HTTP, DNS o TCP se pueden usar en diferentes dispositivos independientemente del sistema operativo con el que cuenta cada uno."

$find = $d.Content.Find
$find.Execute($oldAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $newAnswer, 2)

# --- Change 2 -----------------------------------------------------------
# Insert the answer to "Describir la diferencia entre los modelos de
# protocolo y modelos de referencia." (two body paragraphs plus a trailing
# blank spacer paragraph), mirroring the layout already used for the
# previous answer in the document (List Paragraph style, same indents).

$anchorText = "Describir la diferencia entre los modelos de protocolo y modelos de referencia."
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchorPara = $candidate
        break
    }
}

$anchorPara.Range.InsertParagraphAfter()

# -- new paragraph 1: "Los modelos de protocolo ..." --
$para1 = $d.Paragraphs.Item($anchorPara.Index + 1)
$para1.Style = "List Paragraph"
$para1.Format.LeftIndent = 69.8
$para1.Format.RightIndent = 58.5
$para1.Format.FirstLineIndent = 0
$para1.Format.SpaceBefore = 4.55
$para1.Format.Alignment = 3

$para1Range = $para1.Range
$para1Range.Font.Name = "Cambria"
$para1Range.Font.NameFarEast = "Cambria"
$para1Range.Font.Size = 14
$para1Range.Font.SizeBi = 14
$para1Range.Text = "Los modelos de protocolo son aquellos que proporcionan modelos detallados que coinciden con la estructura de una suite de protocolo en particular. Esto significa mayor especificación y detalle en cuanto a la funcionalidad requerida para interconectar la red humana con la red de datos. Un ejemplo de esto sería el modelo TCP/IP ya que describe las funciones que se producen en cada capa de los protocolos que se encuentran dentro del conjunto TCP/IP."

# -- new paragraph 2: "En cuanto a los modelos de referencia ..." --
$para1.Range.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item($para1.Index + 1)
$para2.Style = "List Paragraph"
$para2.Format.LeftIndent = 69.8
$para2.Format.RightIndent = 58.5
$para2.Format.FirstLineIndent = 0
$para2.Format.SpaceBefore = 4.55
$para2.Format.Alignment = 3

$para2Range = $para2.Range
$para2Range.Font.Name = "Cambria"
$para2Range.Font.NameFarEast = "Cambria"
$para2Range.Font.Size = 14
$para2Range.Font.SizeBi = 14
$para2Range.Text = "En cuanto a los modelos de referencia se tiene que proporcionar un nivel de abstracción mucho menor, generando así una referencia común para mantener consistencia en todos los tipos de protocolos y servicios de la red. Es decir que en estos modelos no se genera tanta minucia ni especificación, por lo tanto, no se puede definir una forma precisa de los servicios de la arquitectura de red. Un ejemplo de esto sería el modelo OSI."

# _GoBack bookmark now sits right at the start of this paragraph's text.
$bmStart = $para2.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmStart))

# -- new paragraph 3: trailing blank spacer paragraph --
$para2.Range.InsertParagraphAfter()
$para3 = $d.Paragraphs.Item($para2.Index + 1)
$para3.Style = "Normal"
$para3.Format.LeftIndent = 69.8
$para3.Format.RightIndent = 58.5
$para3.Format.SpaceBefore = 4.55
$para3.Format.Alignment = 3

$para3Range = $para3.Range
$para3Range.Font.Name = "Cambria"
$para3Range.Font.NameFarEast = "Cambria"
$para3Range.Font.Size = 14
$para3Range.Font.SizeBi = 14

# --- Change 3 -----------------------------------------------------------
# Re-assert the "Práctica trabajo colaborativo" heading text so the stale
# lastRenderedPageBreak marker gets dropped when the run is rebuilt.
$find2 = $d.Content.Find
$find2.Execute("Práctica trabajo colaborativo", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Práctica trabajo colaborativo", 2)
